# Update the "Training Dashboard" sheet with the 04-Nov-2025 progress
# snapshot: PERIOD TO EXPIRE (col H) drops by one day and LAST UPDATE
# (col I) moves from 03-Nov-2025 to 04-Nov-2025 for every data row (3-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Make sure the "LAST UPDATE" column keeps storing literal text (it was
# entered as plain text, not a real date) instead of Excel auto-coercing
# the "dd-MMM-yyyy" looking string into a date serial number.
$ws.Range("I3:I32").NumberFormat = "@"

for ($r = 3; $r -le 32; $r++) {
    $periodCell = $ws.Cells.Item($r, 8)
    $currentPeriod = $periodCell.Value()
    $periodCell.Value = $currentPeriod - 1

    $ws.Cells.Item($r, 9).Value = "04-Nov-2025"
}
